# Apply the diff: update recon/charge columns for rows 34-38 and append
# two new transaction rows (39, 40) to the atom_report_cleaned sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# Rows 34, 35, 37, 38 share the same edit pattern:
#   - E<n> (net_amount_to_be_paid) gets filled with the same value as D<n>
#   - G<n> (settlement_date) gets filled with "12-Nov-2025 00:00:00"
#   - V<n> (recon_status) changes from "NRNS" to "RS"
#   - AA<n> (txn_charges) = 5, AB<n> (gst_18) = 0.9, AE<n> (total_chargeable) = 5.9
# ---------------------------------------------------------------------

$rowsStandard = @(34, 35, 37, 38)
foreach ($r in $rowsStandard) {
    $amount = $ws.Range("D$r").Value()
    $ws.Range("E$r").Value = $amount
    $ws.Range("G$r").Value = "12-Nov-2025 00:00:00"
    $ws.Range("V$r").Value = "RS"
    $ws.Range("AA$r").Value = 5
    $ws.Range("AB$r").Value = 0.9
    $ws.Range("AE$r").Value = 5.9
}

# Row 36 is similar but does NOT get a settlement_date, and its
# recon_status becomes "RNS" (not "RS").
$amount36 = $ws.Range("D36").Value()
$ws.Range("E36").Value = $amount36
$ws.Range("V36").Value = "RNS"
$ws.Range("AA36").Value = 5
$ws.Range("AB36").Value = 0.9
$ws.Range("AE36").Value = 5.9

# ---------------------------------------------------------------------
# New row 39: VELAMALA RUTVIK VIHAAN
# ---------------------------------------------------------------------
$ws.Range("A39").Value = 17337
$ws.Range("B39").Value = "VELAMALA RUTVIK VIHAAN"
$ws.Range("C39").Value = 9603662924
$ws.Range("D39").Value = 6750
$ws.Range("F39").Value = "12-Nov-2025 09:11:42"
$ws.Range("H39").Value = "TRANSACTION IS SUCCESSFUL"
$ws.Range("I39").Value = "OK"
$ws.Range("J39").Value = "PREKGUKG"
$ws.Range("K39").Value = "six thousand seven hundred fifty"
$ws.Range("L39").Value = 100000036600
$ws.Range("M39").Value = "SALESIAN EDUCATION SOCIETY"
$ws.Range("N39").Value = 753702
$ws.Range("O39").Value = 1234
$ws.Range("P39").Value = 11000316275259
$ws.Range("Q39").Value = 1762918783
$ws.Range("R39").Value = 108557239779
$ws.Range("S39").Value = "INR"
$ws.Range("T39").Value = "sale"
$ws.Range("U39").Value = "ICICI UPI QR"
$ws.Range("V39").Value = "NRNS"
$ws.Range("W39").Value = "SIBL0000899"
$ws.Range("X39").Value = "MERCHANT"
$ws.Range("Y39").Value = "UPI"
$ws.Range("Z39").Value = "kotakschoolvsp@gmail.com"
$ws.Range("AC39").Value = 0
$ws.Range("AD39").Value = 0
$ws.Range("AF39").Value = "KOTAK SALESIAN SCHOOL MANAGEMENT ACCOUNT"
$ws.Range("AH39").Value = "REGULAR"
$ws.Range("AI39").Value = 20052
$ws.Range("AJ39").NumberFormat = "@"
$ws.Range("AJ39").Value = "263452"
$ws.Range("AK39").NumberFormat = "@"
$ws.Range("AK39").Value = "2039"

# ---------------------------------------------------------------------
# New row 40: PEDDINTI JASWANTH
# ---------------------------------------------------------------------
$ws.Range("A40").Value = 16211
$ws.Range("B40").Value = "PEDDINTI JASWANTH"
$ws.Range("C40").Value = 9550130840
$ws.Range("D40").Value = 8350
$ws.Range("F40").Value = "12-Nov-2025 10:32:49"
$ws.Range("H40").Value = "TRANSACTION IS SUCCESSFUL"
$ws.Range("I40").Value = "OK"
$ws.Range("J40").Value = "lVl"
$ws.Range("K40").Value = "eight thousand three hundred fifty"
$ws.Range("L40").Value = 100000036600
$ws.Range("M40").Value = "SALESIAN EDUCATION SOCIETY"
$ws.Range("N40").Value = 753702
$ws.Range("O40").Value = 1234
$ws.Range("P40").Value = 11000316293921
$ws.Range("Q40").Value = 1762923744
$ws.Range("R40").Value = 214127184653
$ws.Range("S40").Value = "INR"
$ws.Range("T40").Value = "sale"
$ws.Range("U40").Value = "ICICI UPI QR"
$ws.Range("V40").Value = "NRNS"
$ws.Range("W40").Value = "SIBL0000899"
$ws.Range("X40").Value = "MERCHANT"
$ws.Range("Y40").Value = "UPI"
$ws.Range("Z40").Value = "kotakschoolvsp@gmail.com"
$ws.Range("AC40").Value = 0
$ws.Range("AD40").Value = 0
$ws.Range("AF40").Value = "KOTAK SALESIAN PRIMARY SCHOOL"
$ws.Range("AH40").Value = "REGULAR"
$ws.Range("AI40").Value = 18814
$ws.Range("AJ40").NumberFormat = "@"
$ws.Range("AJ40").Value = "265854"
$ws.Range("AK40").NumberFormat = "@"
$ws.Range("AK40").Value = "2050"
$ws.Range("AL40").Value = "UPI INTENT"
